$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 16:34"

# Row 4
$ws.Cells.Item(4, 2).Value = 4435826
$ws.Cells.Item(4, 3).Value = 2416
$ws.Cells.Item(4, 4).Value = 2137986
$ws.Cells.Item(4, 5).Value = 2147312
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 84
$ws.Cells.Item(4, 8).Value = 150528

# Row 6
$ws.Cells.Item(6, 2).Value = 1516129
$ws.Cells.Item(6, 3).Value = 33626
$ws.Cells.Item(6, 4).Value = 970642
$ws.Cells.Item(6, 5).Value = 511639
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 400
$ws.Cells.Item(6, 8).Value = 33848

# Row 23
$ws.Cells.Item(23, 4).Value = 75083
$ws.Cells.Item(23, 5).Value = 89251
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 23
$ws.Cells.Item(23, 8).Value = 3082

# Row 62
$ws.Cells.Item(62, 2).Value = 24520
$ws.Cells.Item(62, 3).Value = 379
$ws.Cells.Item(62, 4).Value = 14047
$ws.Cells.Item(62, 5).Value = 9922
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 8
$ws.Cells.Item(62, 8).Value = 551

# Row 64
$ws.Cells.Item(64, 2).Value = 21699
$ws.Cells.Item(64, 3).Value = 490
$ws.Cells.Item(64, 4).Value = 12026
$ws.Cells.Item(64, 5).Value = 9549
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 3
$ws.Cells.Item(64, 8).Value = 124

# Row 79
$ws.Cells.Item(79, 2).Value = 11496
$ws.Cells.Item(79, 3).Value = 72
$ws.Cells.Item(79, 4).Value = 6001
$ws.Cells.Item(79, 5).Value = 4770
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 5
$ws.Cells.Item(79, 8).Value = 725

# Row 86
$ws.Cells.Item(86, 2).Value = 9142
$ws.Cells.Item(86, 3).Value = 10
$ws.Cells.Item(86, 4).Value = 8752
$ws.Cells.Item(86, 5).Value = 135
$ws.Cells.Item(86, 6).Value = 0

# Row 92
$ws.Cells.Item(92, 2).Value = 7276
$ws.Cells.Item(92, 3).Value = 41
$ws.Cells.Item(92, 4).Value = 6065
$ws.Cells.Item(92, 5).Value = 1151
$ws.Cells.Item(92, 6).Value = 0

# Row 98
$ws.Cells.Item(98, 1).Value = "Zambia"
$ws.Cells.Item(98, 2).Value = 5002
$ws.Cells.Item(98, 3).Value = 450
$ws.Cells.Item(98, 4).Value = 3195
$ws.Cells.Item(98, 5).Value = 1665
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 142

# Row 99
$ws.Cells.Item(99, 1).Value = "Croacia"
$ws.Cells.Item(99, 2).Value = 4923
$ws.Cells.Item(99, 3).Value = 42
$ws.Cells.Item(99, 4).Value = 4034
$ws.Cells.Item(99, 5).Value = 749
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 140

# Row 100
$ws.Cells.Item(100, 1).Value = "Albania"
$ws.Cells.Item(100, 2).Value = 4880
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 2745
$ws.Cells.Item(100, 5).Value = 1991
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 8).Value = 144

# Row 101
$ws.Cells.Item(101, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(101, 2).Value = 4599
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 1546
$ws.Cells.Item(101, 5).Value = 2994
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 8).Value = 59

# Row 115
$ws.Cells.Item(115, 1).Value = "Hong Kong"
$ws.Cells.Item(115, 2).Value = 2885
$ws.Cells.Item(115, 3).Value = 106
$ws.Cells.Item(115, 4).Value = 1527
$ws.Cells.Item(115, 5).Value = 1335
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 23

# Row 116
$ws.Cells.Item(116, 1).Value = "Libia"
$ws.Cells.Item(116, 2).Value = 2827
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 577
$ws.Cells.Item(116, 5).Value = 2186
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 8).Value = 64

# Row 117
$ws.Cells.Item(117, 1).Value = "Sri Lanka"
$ws.Cells.Item(117, 2).Value = 2807
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 4).Value = 2296
$ws.Cells.Item(117, 5).Value = 500
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 8).Value = 11

# Row 119
$ws.Cells.Item(119, 2).Value = 2555
$ws.Cells.Item(119, 3).Value = 23
$ws.Cells.Item(119, 4).Value = 2352
$ws.Cells.Item(119, 5).Value = 116
$ws.Cells.Item(119, 6).Value = 0

# Row 129
$ws.Cells.Item(129, 1).Value = "Namibia"
$ws.Cells.Item(129, 2).Value = 1917
$ws.Cells.Item(129, 3).Value = 74
$ws.Cells.Item(129, 4).Value = 104
$ws.Cells.Item(129, 5).Value = 1805
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 8).Value = 8

# Row 130
$ws.Cells.Item(130, 1).Value = "Ruanda"
$ws.Cells.Item(130, 2).Value = 1879
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 975
$ws.Cells.Item(130, 5).Value = 899
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 8).Value = 5

# Row 131
$ws.Cells.Item(131, 1).Value = "Islandia"
$ws.Cells.Item(131, 2).Value = 1857
$ws.Cells.Item(131, 3).Value = 3
$ws.Cells.Item(131, 4).Value = 1823
$ws.Cells.Item(131, 5).Value = 24
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 8).Value = 10

# Row 168
$ws.Cells.Item(168, 4).Value = 293
$ws.Cells.Item(168, 5).Value = 51
$ws.Cells.Item(168, 6).Value = 0
